$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1063.083459787557
$ws.Range("D2").Value = 0.9292278968751633
$ws.Range("E2").Value = 2.102326757713717
$ws.Range("F2").Value = 1.57970949319211
$ws.Range("G2").Value = -0.004993195227396457
$ws.Range("H2").Value = -0.01331228679955836
$ws.Range("I2").Value = -0.634720221722993
$ws.Range("J2").Value = 1.22193563099095
$ws.Range("K2").Value = [double]"-1.60549930007693e-15"
$ws.Range("L2").Value = [double]"-8.449207163177904e-15"
$ws.Range("M2").Value = 0.3836966487215132
$ws.Range("N2").Value = [double]"3.829904414741309e-16"
$ws.Range("O2").Value = 40.4243803743045
$ws.Range("P2").Value = 362.0371235812341
$ws.Range("Q2").Value = 361.7735046098999
$ws.Range("S2").Value = 0.4344036265094725
$ws.Range("T2").Value = 0.05192932268742551
$ws.Range("U2").Value = 0.5986000988246902
$ws.Range("V2").Value = -1.57970949319211
$ws.Range("W2").Value = 1.615307186802609
$ws.Range("X2").Value = 1977
$ws.Range("Z2").Value = 0.3566017316017316
